# Auto-generated edit script applying numeric value changes per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# hunk1: ALC row 41
$ws.Cells.Item(41, 8).Value = 314.66666
$ws.Cells.Item(41, 9).Value = 368.75
$ws.Cells.Item(41, 10).Value = 252.85715
$ws.Cells.Item(41, 11).Value = 368.75
$ws.Cells.Item(41, 12).Value = 252.85715
$ws.Cells.Item(41, 13).Value = 71.25
$ws.Cells.Item(41, 14).Value = -1132.85715

# hunk2: ALC row 53
$ws.Cells.Item(53, 8).Value = 998.3
$ws.Cells.Item(53, 9).Value = 1284.9333
$ws.Cells.Item(53, 11).Value = 1284.9333
$ws.Cells.Item(53, 13).Value = -647.9332999999999

# hunk3: ALC row 110
$ws.Cells.Item(110, 8).Value = 41242.855
$ws.Cells.Item(110, 10).Value = 41242.855
$ws.Cells.Item(110, 12).Value = 41242.855
$ws.Cells.Item(110, 14).Value = -49422.855

# hunk4: ALC row 114
$ws.Cells.Item(114, 8).Value = 0
$ws.Cells.Item(114, 10).Value = 0
$ws.Cells.Item(114, 12).Value = 0
$ws.Cells.Item(114, 14).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# hunk5: ARM row 88
$ws.Cells.Item(88, 8).Value = 2382.5386
$ws.Cells.Item(88, 9).Value = 2361
$ws.Cells.Item(88, 10).Value = 2401
$ws.Cells.Item(88, 11).Value = 2361
$ws.Cells.Item(88, 12).Value = 2401
$ws.Cells.Item(88, 13).Value = -1955
$ws.Cells.Item(88, 14).Value = -3213

# hunk6: ARM row 91
$ws.Cells.Item(91, 8).Value = 2382.5386
$ws.Cells.Item(91, 9).Value = 2361
$ws.Cells.Item(91, 10).Value = 2401
$ws.Cells.Item(91, 11).Value = 2361
$ws.Cells.Item(91, 12).Value = 2401
$ws.Cells.Item(91, 13).Value = -957
$ws.Cells.Item(91, 14).Value = -5209

# hunk7: ARM row 102
$ws.Cells.Item(102, 8).Value = 1495.56
$ws.Cells.Item(102, 9).Value = 1286.1875
$ws.Cells.Item(102, 11).Value = 1286.1875
$ws.Cells.Item(102, 13).Value = 335.8125

# hunk8: ARM row 113
$ws.Cells.Item(113, 8).Value = 30057.143
$ws.Cells.Item(113, 10).Value = 30057.143
$ws.Cells.Item(113, 12).Value = 30057.143
$ws.Cells.Item(113, 14).Value = -38735.143

# hunk9: ARM row 114
$ws.Cells.Item(114, 8).Value = 30000
$ws.Cells.Item(114, 10).Value = 30000
$ws.Cells.Item(114, 12).Value = 30000
$ws.Cells.Item(114, 14).Value = -38678

# hunk10: ARM row 117
$ws.Cells.Item(117, 8).Value = 29333.334
$ws.Cells.Item(117, 10).Value = 29333.334
$ws.Cells.Item(117, 12).Value = 29333.334
$ws.Cells.Item(117, 14).Value = -38511.334

# hunk11: ARM row 119
$ws.Cells.Item(119, 8).Value = 31200
$ws.Cells.Item(119, 10).Value = 31200
$ws.Cells.Item(119, 12).Value = 31200
$ws.Cells.Item(119, 14).Value = -40876

$ws = $wb.Worksheets.Item("BSM")
# hunk12: BSM row 99
$ws.Cells.Item(99, 8).Value = 1796.138
$ws.Cells.Item(99, 9).Value = 1625.7368
$ws.Cells.Item(99, 10).Value = 2119.9
$ws.Cells.Item(99, 11).Value = 1625.7368
$ws.Cells.Item(99, 12).Value = 2119.9
$ws.Cells.Item(99, 13).Value = -127.7367999999999
$ws.Cells.Item(99, 14).Value = -5115.9

# hunk13: BSM row 105
$ws.Cells.Item(105, 8).Value = 2300.7925
$ws.Cells.Item(105, 9).Value = 2264.2856
$ws.Cells.Item(105, 10).Value = 2324.75
$ws.Cells.Item(105, 11).Value = 2264.2856
$ws.Cells.Item(105, 12).Value = 2324.75
$ws.Cells.Item(105, 13).Value = -517.2856000000002
$ws.Cells.Item(105, 14).Value = -5818.75

# hunk14: BSM row 134
$ws.Cells.Item(134, 8).Value = 574555.1
$ws.Cells.Item(134, 9).Value = 977899.9
$ws.Cells.Item(134, 10).Value = 4309.069
$ws.Cells.Item(134, 11).Value = 2933699.7
$ws.Cells.Item(134, 12).Value = 12927.207
$ws.Cells.Item(134, 13).Value = -2931164.7
$ws.Cells.Item(134, 14).Value = -17997.207

$ws = $wb.Worksheets.Item("CRP")
# hunk15: CRP row 10
$ws.Cells.Item(10, 8).Value = 16902.334
$ws.Cells.Item(10, 9).Value = 353.5
$ws.Cells.Item(10, 10).Value = 50000
$ws.Cells.Item(10, 11).Value = 353.5
$ws.Cells.Item(10, 12).Value = 50000
$ws.Cells.Item(10, 13).Value = -214.5
$ws.Cells.Item(10, 14).Value = -50278

# hunk16: CRP row 58
$ws.Cells.Item(58, 8).Value = 2748.7966
$ws.Cells.Item(58, 9).Value = 3186.0698
$ws.Cells.Item(58, 10).Value = 1573.625
$ws.Cells.Item(58, 11).Value = 3186.0698
$ws.Cells.Item(58, 12).Value = 1573.625
$ws.Cells.Item(58, 13).Value = -2983.0698
$ws.Cells.Item(58, 14).Value = -1979.625

# hunk17: CRP row 107
$ws.Cells.Item(107, 8).Value = 1927.826
$ws.Cells.Item(107, 9).Value = 662.8
$ws.Cells.Item(107, 10).Value = 4299.75
$ws.Cells.Item(107, 11).Value = 662.8
$ws.Cells.Item(107, 12).Value = 4299.75
$ws.Cells.Item(107, 13).Value = 1257.2
$ws.Cells.Item(107, 14).Value = -8139.75

# hunk18: CRP row 127
$ws.Cells.Item(127, 8).Value = 31884.285
$ws.Cells.Item(127, 10).Value = 31884.285
$ws.Cells.Item(127, 12).Value = 31884.285
$ws.Cells.Item(127, 14).Value = -41804.285

# hunk19: CRP row 134
$ws.Cells.Item(134, 8).Value = 17188970
$ws.Cells.Item(134, 9).Value = 1853240.4
$ws.Cells.Item(134, 11).Value = 5559721.199999999
$ws.Cells.Item(134, 13).Value = -5557186.199999999

# hunk20: CRP row 136
$ws.Cells.Item(136, 8).Value = 2748.7966
$ws.Cells.Item(136, 9).Value = 3186.0698
$ws.Cells.Item(136, 10).Value = 1573.625
$ws.Cells.Item(136, 11).Value = 9558.2094
$ws.Cells.Item(136, 12).Value = 4720.875
$ws.Cells.Item(136, 13).Value = -7008.2094
$ws.Cells.Item(136, 14).Value = -9820.875

$ws = $wb.Worksheets.Item("CUL")
# hunk21: CUL row 123
$ws.Cells.Item(123, 8).Value = 1955
$ws.Cells.Item(123, 9).Value = 966
$ws.Cells.Item(123, 10).Value = 6900
$ws.Cells.Item(123, 11).Value = 2898
$ws.Cells.Item(123, 12).Value = 20700
$ws.Cells.Item(123, 13).Value = -448
$ws.Cells.Item(123, 14).Value = -25600

# hunk22: CUL row 131
$ws.Cells.Item(131, 8).Value = 907.23
$ws.Cells.Item(131, 10).Value = 916.2283
$ws.Cells.Item(131, 12).Value = 2748.6849
$ws.Cells.Item(131, 14).Value = -12828.6849

# hunk23: CUL row 134
$ws.Cells.Item(134, 8).Value = 6403.2563
$ws.Cells.Item(134, 9).Value = 5086.0557
$ws.Cells.Item(134, 10).Value = 7532.2856
$ws.Cells.Item(134, 11).Value = 15258.1671
$ws.Cells.Item(134, 12).Value = 22596.8568
$ws.Cells.Item(134, 13).Value = -10188.1671
$ws.Cells.Item(134, 14).Value = -32736.8568

$ws = $wb.Worksheets.Item("GSM")
# hunk24: GSM row 7
$ws.Cells.Item(7, 8).Value = 3001
$ws.Cells.Item(7, 9).Value = 3001
$ws.Cells.Item(7, 11).Value = 3001
$ws.Cells.Item(7, 13).Value = -2889

# hunk25: GSM row 8
$ws.Cells.Item(8, 8).Value = 3001
$ws.Cells.Item(8, 9).Value = 3001
$ws.Cells.Item(8, 11).Value = 3001
$ws.Cells.Item(8, 13).Value = -2862

# hunk26: GSM row 80
$ws.Cells.Item(80, 8).Value = 2309.2122
$ws.Cells.Item(80, 9).Value = 2280.3914
$ws.Cells.Item(80, 10).Value = 2375.5
$ws.Cells.Item(80, 11).Value = 2280.3914
$ws.Cells.Item(80, 12).Value = 2375.5
$ws.Cells.Item(80, 13).Value = -1282.3914
$ws.Cells.Item(80, 14).Value = -4371.5

# hunk27: GSM row 83
$ws.Cells.Item(83, 8).Value = 2309.2122
$ws.Cells.Item(83, 9).Value = 2280.3914
$ws.Cells.Item(83, 10).Value = 2375.5
$ws.Cells.Item(83, 11).Value = 11401.957
$ws.Cells.Item(83, 12).Value = 11877.5
$ws.Cells.Item(83, 13).Value = -6409.957
$ws.Cells.Item(83, 14).Value = -21861.5

# hunk28: GSM row 110
$ws.Cells.Item(110, 8).Value = 32800
$ws.Cells.Item(110, 10).Value = 32800
$ws.Cells.Item(110, 12).Value = 32800
$ws.Cells.Item(110, 14).Value = -40980

# hunk29: GSM row 116
$ws.Cells.Item(116, 8).Value = 0
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 14).ClearContents()

# hunk30: GSM row 120
$ws.Cells.Item(120, 8).Value = 0
$ws.Cells.Item(120, 10).Value = 0
$ws.Cells.Item(120, 12).Value = 0
$ws.Cells.Item(120, 14).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# hunk31: LTW row 2
$ws.Cells.Item(2, 8).Value = 108352.94
$ws.Cells.Item(2, 10).Value = 28000
$ws.Cells.Item(2, 12).Value = 28000
$ws.Cells.Item(2, 14).Value = -28224

# hunk32: LTW row 46
$ws.Cells.Item(46, 8).Value = 821.96155
$ws.Cells.Item(46, 9).Value = 1205.1666
$ws.Cells.Item(46, 10).Value = 707
$ws.Cells.Item(46, 11).Value = 1205.1666
$ws.Cells.Item(46, 12).Value = 707
$ws.Cells.Item(46, 13).Value = -1017.1666
$ws.Cells.Item(46, 14).Value = -1083

# hunk33: LTW row 132
$ws.Cells.Item(132, 8).Value = 4896.864
$ws.Cells.Item(132, 9).Value = 8556
$ws.Cells.Item(132, 10).Value = 2363.6155
$ws.Cells.Item(132, 11).Value = 25668
$ws.Cells.Item(132, 12).Value = 7090.8465
$ws.Cells.Item(132, 13).Value = -23138
$ws.Cells.Item(132, 14).Value = -12150.8465

# hunk34: LTW row 136
$ws.Cells.Item(136, 8).Value = 1528.8823
$ws.Cells.Item(136, 9).Value = 1028.3478
$ws.Cells.Item(136, 10).Value = 2575.4546
$ws.Cells.Item(136, 11).Value = 3085.0434
$ws.Cells.Item(136, 12).Value = 7726.3638
$ws.Cells.Item(136, 13).Value = -535.0434
$ws.Cells.Item(136, 14).Value = -12826.3638

$ws = $wb.Worksheets.Item("WVR")
# hunk35: WVR row 9
$ws.Cells.Item(9, 8).Value = 1883.3334
$ws.Cells.Item(9, 9).Value = 1900
$ws.Cells.Item(9, 10).Value = 1800
$ws.Cells.Item(9, 11).Value = 1900
$ws.Cells.Item(9, 12).Value = 1800
$ws.Cells.Item(9, 13).Value = -1760
$ws.Cells.Item(9, 14).Value = -2080

# hunk36: WVR row 132
$ws.Cells.Item(132, 8).Value = 2920.5366
$ws.Cells.Item(132, 9).Value = 3497.56
$ws.Cells.Item(132, 10).Value = 2018.9375
$ws.Cells.Item(132, 11).Value = 10492.68
$ws.Cells.Item(132, 12).Value = 6056.8125
$ws.Cells.Item(132, 13).Value = -7962.68
$ws.Cells.Item(132, 14).Value = -11116.8125

# hunk37: WVR row 136
$ws.Cells.Item(136, 8).Value = 1231.2982
$ws.Cells.Item(136, 9).Value = 862
$ws.Cells.Item(136, 10).Value = 1739.0834
$ws.Cells.Item(136, 11).Value = 2586
$ws.Cells.Item(136, 12).Value = 5217.2502
$ws.Cells.Item(136, 13).Value = -36
$ws.Cells.Item(136, 14).Value = -10317.2502
